$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.855.44"
$ws.Range("E2").Value = "  -1.92%  "
$ws.Range("D3").Value = "1.824.86"
$ws.Range("E3").Value = "  -2.20%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.9998"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "239.16"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.70%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.6899"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.08%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.9999"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.13%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.07600"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -3.21%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.3013"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -3.89%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "23.29"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -4.74%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07710"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -3.85%  "
$ws.Range("D12").Value = "1.826.38"
$ws.Range("E12").Value = "  -3.06%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "5.030"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -3.13%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "89.76"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -3.70%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.6699"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -4.28%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "6.365"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.25%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.000008256"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").Value = "28.859.27"
$ws.Range("E18").Value = "  -2.05%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "241.94"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -4.41%  "
$ws.Range("D20").Value = "2.079.20"
$ws.Range("E20").Value = "  -2.66%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "12.57"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -4.27%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.9993"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.19%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "7.363"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -3.06%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.9998"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.15%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.1466"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -5.66%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "160.43"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.14%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "8.692"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -3.55%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "18.11"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -3.38%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.527"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.81%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.173"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -3.46%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.118"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -3.46%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.193"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.69%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.05084"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -4.02%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.7424"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.72%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.804"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -4.15%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.136"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.42%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.684"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.25%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.01830"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -2.39%  "
$ws.Range("D39").Value = "1.193.04"
$ws.Range("E39").Value = "  -5.20%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.669"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -2.76%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.9130"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.65%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "108.04"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.54%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.9993"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("D44").Value = "1.978.35"
$ws.Range("E44").Value = "  -2.92%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.5161"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("E46").Value = "  -5.29%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "9.424"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.73%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "5.200"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -12.65%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.720"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -3.76%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "62.02"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -13.13%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.4180"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -2.82%  "
